$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

# Row 11 / column B ("Rule" name) changes from "R40" to the text "1".
# A direct Value assignment of a numeric-looking string like "1" would be
# auto-coerced into a genuine number by the type-inference used for
# Range.Value, which would store the cell as a number (and would also
# require a new/augmented style to force text via quote-prefix). To end
# up with a true *text* cell - matching the original "string" column
# semantics (t="s") and keeping the existing cell style - compute the
# text through a formula and commit it as a literal value via a
# copy/paste-special (values only), which bypasses the COM marshalling
# round trip that would otherwise re-interpret "1" as numeric.
$ws.Range("B11").Formula = "=TEXT(1,""0"")"
$ws.Calculate()
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$wb.Save()
